$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: copy a single cell (value + formatting) from one location to
# another.  Copying cell-by-cell (instead of whole merged ranges at once)
# avoids a runtime quirk that duplicates style entries when a merged range
# is copied as a block.
# ---------------------------------------------------------------------------
function Copy-Cell($srcRow, $srcCol, $dstRow, $dstCol) {
    $ws.Cells.Item($srcRow, $srcCol).Copy($ws.Cells.Item($dstRow, $dstCol))
}

function Copy-RowCells($srcRow, $dstRow, $maxCol) {
    for ($c = 1; $c -le $maxCol; $c++) {
        Copy-Cell $srcRow $c $dstRow $c
    }
}

# Set a text value into a cell while preserving its (possibly numeric)
# number format: temporarily switch the format to Text ("@"), assign the
# value (so it is stored as a shared string instead of being coerced into
# a number), then restore the original number format.
function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $origFmt = $cell.NumberFormat()
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = $origFmt
}

# ---------------------------------------------------------------------------
# Step 1: relocate the existing "totals" row (old row 9) and "footer" row
# (old row 10) down to their new positions (rows 13 and 14) to make room for
# four new item rows.
# ---------------------------------------------------------------------------
Copy-RowCells 10 14 17
Copy-RowCells 9 13 17

# ---------------------------------------------------------------------------
# Step 2: build the four new item rows (9-12), using row 8 as the
# formatting template (same column layout/merges as every other item row).
# ---------------------------------------------------------------------------
Copy-RowCells 8 9 17
Copy-RowCells 8 10 17
Copy-RowCells 8 11 17
Copy-RowCells 8 12 17

# Re-create the row-level merges for the four new rows (mirrors rows 7/8).
$ws.Range("A9:B9").Merge()
$ws.Range("C9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("N9:O9").Merge()

$ws.Range("A10:B10").Merge()
$ws.Range("C10:G10").Merge()
$ws.Range("H10:K10").Merge()
$ws.Range("L10:M10").Merge()
$ws.Range("N10:O10").Merge()

$ws.Range("A11:B11").Merge()
$ws.Range("C11:G11").Merge()
$ws.Range("H11:K11").Merge()
$ws.Range("L11:M11").Merge()
$ws.Range("N11:O11").Merge()

$ws.Range("A12:B12").Merge()
$ws.Range("C12:G12").Merge()
$ws.Range("H12:K12").Merge()
$ws.Range("L12:M12").Merge()
$ws.Range("N12:O12").Merge()

# ---------------------------------------------------------------------------
# Step 3: fill in the actual item data for the 4 new rows + update row 8's
# item (it keeps its slot, but its row index number stays "2" already).
# ---------------------------------------------------------------------------

# Row 8 -> CONCOR PLUS 10/25MG 30 F.C. TABLETS (was CONTAFEVER; that item now
# moves down to row 9)
$ws.Cells.Item(8,1).Value = 2
Set-TextValue 8 3  "CONCOR PLUS 10/25MG 30 F.C. TABLETS"
Set-TextValue 8 8  "2:1"
Set-TextValue 8 12 "1"
Set-TextValue 8 14 "108.00"
Set-TextValue 8 16 "35.6400"
Set-TextValue 8 17 "0:1"

# Row 9 -> CONTAFEVER N 200MG/5ML SUSP. 120ML
$ws.Cells.Item(9,1).Value = 3
Set-TextValue 9 3  "CONTAFEVER N 200MG/5ML SUSP. 120ML"
Set-TextValue 9 8  "11:0"
Set-TextValue 9 12 "1"
Set-TextValue 9 14 "33.00"
Set-TextValue 9 16 "33.0000"
Set-TextValue 9 17 "1:0"

# Row 10 -> OMEZ 20MG 14 CAPS.
$ws.Cells.Item(10,1).Value = 4
Set-TextValue 10 3  "OMEZ 20MG 14 CAPS."
Set-TextValue 10 8  "1:1"
Set-TextValue 10 12 "1"
Set-TextValue 10 14 "56.00"
Set-TextValue 10 16 "28.0000"
Set-TextValue 10 17 "0:1"

# Row 11 -> (belt) حزام فقرات
$ws.Cells.Item(11,1).Value = 5
Set-TextValue 11 3  "حزام فقرات "
Set-TextValue 11 8  "0:0"
Set-TextValue 11 12 "0"
Set-TextValue 11 14 "130.00"
Set-TextValue 11 16 "130.0000"
Set-TextValue 11 17 "1:0"

# Row 12 -> (extra long socks) صوفي طويل جدا
$ws.Cells.Item(12,1).Value = 6
Set-TextValue 12 3  "صوفي طويل جدا"
Set-TextValue 12 8  "14:0"
Set-TextValue 12 12 "0"
Set-TextValue 12 14 "50.00"
Set-TextValue 12 16 "50.0000"
Set-TextValue 12 17 "1:0"

# ---------------------------------------------------------------------------
# Step 4: update the relocated totals row (now row 13) with the new sum of
# the "selling price" column, and the relocated footer row (now row 14)
# with the refreshed timestamp.
# ---------------------------------------------------------------------------
$ws.Cells.Item(13,16).Value = 345.64

Set-TextValue 14 1 "Saturday, 4 October, 2025 10:56 AM"

Write-Host "Edit complete"
